# 1.0.9 - Boksi tekstuurit
# Adds a new "Testi 6" block (rows 68-78) to the testing-log worksheet,
# mirroring the existing "Testi 5" block's layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reproduce the formatting (fills/borders/number format) of the previous
# test block (rows 56-66) onto the new block (rows 68-78) by copy/paste
# of formats only - this keeps the underlying cell styles identical to
# the existing pattern used throughout the sheet.
$ws.Range("A56:B66").Copy() | Out-Null
$ws.Range("A68:B78").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Section header row
$ws.Range("A68").Value = "Testi 6"

# Versio
$ws.Range("A69").Value = "Versio"
$ws.Range("B69").Value = "1.0.9"

# Päivämäärä (16.9.2020)
$ws.Range("A70").Value = "Päivämäärä"
$ws.Range("B70").Value = 44090

# Testaaja
$ws.Range("A71").Value = "Testaaja"
$ws.Range("B71").Value = "Jasper"

# Testaustyyppi
$ws.Range("A72").Value = "Testaustyyppi"
$ws.Range("B72").Value = "Toiminnallisuustesti"

# Mitä testataan?
$ws.Range("A73").Value = "Mitä testataan?"
$ws.Range("B73").Value = "Boksien tekstuureita ja niiden istuvuutta"

# Odotetut tulokset
$ws.Range("A74").Value = "Odotetut tulokset"
$ws.Range("B74").Value = "Boksien tekstuurit toimivat normaalisti"

# Testin tulos
$ws.Range("A75").Value = "Testin tulos"
$ws.Range("B75").Value = "Boksien tekstuurit toimivat  "

# Toimenpiteet
$ws.Range("A76").Value = "Toimenpiteet"
$ws.Range("B76").Value = "N/A"

# Testausympäristö
$ws.Range("A77").Value = "Testausympäristö"
$ws.Range("B77").Value = "Unity"

# Kommentit
$ws.Range("A78").Value = "Kommentit"
$ws.Range("B78").Value = "Lisättiin bokseille uusi ulkonäkö onnistuneesti"

# Match the saved view state: scrolled down with B78 selected.
$ws.Range("B78").Select() | Out-Null
